$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-15 21:00:55"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-15 21:00:50"
$wsZhCn.Range("K2").Value = "2016-08-15 21:01:17"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-15 21:00:55"
$wsDeDe.Range("K2").Value = "2016-08-15 21:01:24"
